$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(PriceText, VolumeText)
$updates = @{
    2 = @('66.258.64', '  +0.36%  ')
    3 = @('3.573.00', '  +2.57%  ')
    4 = @('1.00', '  +0.14%  ')
    5 = @('608.91', '  +0.82%  ')
    6 = @('145.29', '  +1.40%  ')
    7 = @('3.571.86', '  +2.50%  ')
    8 = @('1.00', '  +0.25%  ')
    9 = @('0.493', '  +3.90%  ')
    10 = @('0.137', '  +1.29%  ')
    11 = @('7.90', '  -3.38%  ')
    12 = @('0.414', '  +0.53%  ')
    13 = @('4.178.22', '  +2.67%  ')
    14 = @('0.0000208', '  +2.37%  ')
    15 = @('29.97', '  -1.49%  ')
    16 = @('3.561.05', '  +2.26%  ')
    17 = @('66.352.09', '  +0.34%  ')
    18 = @('0.115', '  -0.92%  ')
    19 = @('11.49', '  +10.74%  ')
    20 = @('6.23', '  +0.93%  ')
    21 = @('14.88', '  +0.99%  ')
    22 = @('429.95', '  +2.28%  ')
    23 = @('0.617', '  +4.36%  ')
    24 = @('79.23', '  +2.31%  ')
    25 = @('3.714.00', '  +2.86%  ')
    26 = @('1.00', '  +0.00%  ')
    27 = @('0.0000118', '  +3.58%  ')
    28 = @('2.51', '  +2.13%  ')
    29 = @('7.94', '  -0.72%  ')
    30 = @('9.09', '  -2.47%  ')
    31 = @('1.00', '  +0.00%  ')
    32 = @('25.67', '  +2.02%  ')
    33 = @('1.46', '  -1.68%  ')
    34 = @('3.565.87', '  +2.53%  ')
    35 = @('0.153', '  -5.70%  ')
    36 = @('1.00', '  +0.06%  ')
    37 = @('1.74', '  +1.37%  ')
    38 = @('7.87', '  +2.15%  ')
    39 = @('5.62', '  +0.56%  ')
    40 = @('177.75', '  +4.53%  ')
    41 = @('0.999', '  +0.12%  ')
    42 = @('0.0849', '  -1.83%  ')
    43 = @('5.23', '  +2.59%  ')
    44 = @('0.898', '  +0.90%  ')
    45 = @('1.93', '  +0.70%  ')
    46 = @('46.19', '  +2.56%  ')
    47 = @('1.21', '  +1.08%  ')
    48 = @('25.73', '  -1.19%  ')
    49 = @('2.40', '  +2.16%  ')
    50 = @('23.59', '  +9.25%  ')
    51 = @('7.15', '  +0.28%  ')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dCell = $ws.Range("D$row")
    $dCell.NumberFormat = "@"   # force text so values like "1.00" keep their literal form
    $dCell.Value = $vals[0]
    $dCell.Style = "Normal"     # drop the temporary text-format style again
    $ws.Range("E$row").Value = $vals[1]
}
